$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: citation could not be resolved -> reset to "unknown"/"not found" placeholders
$ws.Range("C2").Value = "Unknown Title"
$ws.Range("D2").Value = "Unknown Abstract"
$ws.Range("E2").Value = "[]"
$ws.Range("F2").Value = "not found"
$ws.Range("G2").Value = "N/A"
# H2 needs the literal text "1970-01-01" (not an Excel date serial). Copying the
# format+value from H5 (which already stores that same text) keeps it a plain
# shared-string cell with no extra number-format style instead of auto-converting
# the typed string into a date.
$ws.Range("H5").Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("J2").Value = ""

# Row 3: fill in the full-text journal info + de-duplicated (single-space) author list
$ws.Range("D3").Value = "Supplemental Digital Content is Available in the Text.`n"
$ws.Range("E3").Value = "[Kenneth A.%Egol%NULL%0, Sanjit R.%Konda%NULL%0, Mackenzie L.%Bird%NULL%0, Nicket%Dedhia%NULL%0, Emma K.%Landes%NULL%0, Rachel A.%Ranson%NULL%0, Sara J.%Solasz%NULL%0, Vinay K.%Aggarwal%NULL%0, Joseph A.%Bosco%NULL%0, David L.%Furgiuele%NULL%0, Abhishek%Ganta%NULL%0, Jason%Gould%NULL%0, Thomas R.%Lyon%NULL%0, Toni M.%McLaurin%NULL%0, Nirmal C.%Tejwani%NULL%0, Joseph D.%Zuckerman%NULL%0, Philipp%Leucht%NULL%0]"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "Journal of Orthopaedic Trauma"

# Row 4: same fix pattern
$ws.Range("D4").Value = "Supplemental Digital Content is Available in the Text.`n"
$ws.Range("E4").Value = "[Drake G.%LeBrun%NULL%0, Maxwell A.%Konnaris%NULL%0, Gregory C.%Ghahramani%NULL%0, Ajay%Premkumar%NULL%0, Chris J.%DeFrancesco%NULL%0, Jordan A.%Gruskay%NULL%0, Aleksey%Dvorzhinskiy%NULL%0, Milan S.%Sandhu%NULL%0, Elan M.%Goldwyn%NULL%0, Christopher L.%Mendias%NULL%0, William M.%Ricci%NULL%0]"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "Journal of Orthopaedic Trauma"

# Row 5: same fix pattern
$ws.Range("D5").Value = "Supplemental Digital Content is Available in the Text.`n"
$ws.Range("E5").Value = "[Amit%Thakrar%NULL%0, Karen%Chui%NULL%0, Akhil%Kapoor%NULL%1, John%Hambidge%NULL%1]"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = "Journal of Orthopaedic Trauma"

# Row 6: citation could not be resolved -> reset to "unknown"/"not found" placeholders
$ws.Range("C6").Value = "Unknown Title"
$ws.Range("D6").Value = "Unknown Abstract"
$ws.Range("E6").Value = "[]"
$ws.Range("F6").Value = "not found"
$ws.Range("G6").Value = "N/A"
$ws.Range("H5").Copy()
$ws.Range("H6").PasteSpecial(-4163)
$ws.Range("J6").Value = ""

# Row 7: de-duplicated author list + publisher info
$ws.Range("E7").Value = "[Karen%Chui%NULL%0, Amit%Thakrar%NULL%0, Shivakumar%Shankar%NULL%1]"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = "The British Editorial Society of Bone and Joint Surgery"
